# Add/update metadata report row for Akurana (JAN 2025) on row 3.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Core descriptive columns
$ws.Range("A3").Value = 2025
$ws.Range("B3").Value = "JAN"
$ws.Range("C3").Value = "10/01-01/01"
$ws.Range("D3").Value = "11/01, 12/01, 13/01, 14/01, 15/01, 16/01, 17/01, 18/01, 19/01, 20/01, 21/01, 22/01, 23/01, 24/01, 25/01, 26/01, 27/01, 28/01, 29/01, 30/01, 31/01"
$ws.Range("E3").Value = "Akurana"

# Remaining metric columns (F through AO) are all "-" for this row.
for ($col = 6; $col -le 41; $col++) {
    $ws.Cells.Item(3, $col).Value = "-"
}
